$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
$win.ScrollRow = 31
Write-Host "ScrollRow final:" $win.ScrollRow
